$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New maintenance comment rows to append starting at row 595 (Date, Comment)
$rows = @(
    @{ Row = 595; Date = [DateTime]"2019-07-17"; Comment = "Monthly sampling" },
    @{ Row = 596; Date = [DateTime]"2019-07-26"; Comment = "Cleaned sensors" },
    @{ Row = 597; Date = [DateTime]"2019-07-30"; Comment = "Calibrated sonde, downloaded data" },
    @{ Row = 598; Date = [DateTime]"2019-08-20"; Comment = "Monthly sampling" },
    @{ Row = 599; Date = [DateTime]"2019-08-28"; Comment = "MR service visit - logger out for the week" },
    @{ Row = 600; Date = [DateTime]"2019-09-19"; Comment = "Montly sampling and cleaned sensors" },
    @{ Row = 601; Date = [DateTime]"2019-09-25"; Comment = "Cleaned sensors" },
    @{ Row = 602; Date = [DateTime]"2019-10-30"; Comment = "Logger our - service visit for MR. Underwater PFD out for checks" },
    @{ Row = 603; Date = [DateTime]"2019-11-05"; Comment = "Logger back out, underwater PFD back out, PRT chain reading dodgy until 8/11/2019" },
    @{ Row = 604; Date = [DateTime]"2019-11-08"; Comment = "PRTr eading again OK, but BP Now reading dodgy" },
    @{ Row = 605; Date = [DateTime]"2019-11-19"; Comment = "Circuit board busted. Sent back to MR" },
    @{ Row = 606; Date = [DateTime]"2020-02-05"; Comment = "PRT box out and sent to MR for replacement" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Comment
}

$ws.Range("B606").Select()
